$d = $word.ActiveDocument

$updates = @(
    ,@(3, "[[PERSON_1]] – „s [[PERSON_1]]“, „o [[PERSON_1]]“")
    ,@(4, "[[PERSON_2]] – „k Evě Marečkové“, „u [[PERSON_2]]“")
    ,@(5, "[[PERSON_3]] – „pro [[PERSON_4]]“, „s [[PERSON_5]]“")
    ,@(6, "[[PERSON_6]] – „s [[PERSON_6]]“, „o [[PERSON_6]]“")
    ,@(7, "[[PERSON_7]] – „u [[PERSON_7]]“, „k [[PERSON_7]]“")
    ,@(8, "[[PERSON_8]] – „o [[PERSON_8]]“, „se [[PERSON_8]]“")
    ,@(9, "[[PERSON_9]] – „k [[PERSON_9]]“, „u [[PERSON_9]]“")
    ,@(10, "[[PERSON_10]] – „s [[PERSON_10]]“, „o [[PERSON_10]]“")
    ,@(11, "[[PERSON_11]] – „u [[PERSON_11]]“, „s [[PERSON_11]]“")
    ,@(12, "[[PERSON_12]] – „s [[PERSON_12]]“, „k [[PERSON_12]]“")
    ,@(13, "[[PERSON_13]] – „s [[PERSON_14]]“, „o [[PERSON_15]]“")
    ,@(14, "[[PERSON_16]] – „ke [[PERSON_16]]“, „o [[PERSON_16]]“")
    ,@(15, "[[PERSON_17]] – „o [[PERSON_17]]“, „s [[PERSON_17]]“")
    ,@(16, "[[PERSON_18]] – „u [[PERSON_18]]“, „s [[PERSON_18]]“")
    ,@(17, "[[PERSON_19]] – „ke [[PERSON_19]]“, „o [[PERSON_19]]“")
    ,@(18, "[[PERSON_20]] – „s [[PERSON_21]]“, „o [[PERSON_20]]“")
    ,@(19, "[[PERSON_22]] – „s [[PERSON_22]]“, „o [[PERSON_22]]“")
    ,@(20, "[[PERSON_23]] – „k [[PERSON_23]]“, „od [[PERSON_23]]“")
    ,@(21, "[[PERSON_24]] – „o [[PERSON_24]]“, „s [[PERSON_24]]“")
    ,@(22, "[[PERSON_25]] – „o [[PERSON_26]]“, „se [[PERSON_27]]“")
    ,@(23, "[[PERSON_28]] – „s [[PERSON_28]]“, „u [[PERSON_28]]“")
    ,@(24, "[[PERSON_29]] – „o [[PERSON_30]]“, „s [[PERSON_31]]“")
    ,@(25, "[[PERSON_32]] – „k [[PERSON_32]]“, „o [[PERSON_32]]“")
    ,@(26, "[[PERSON_33]] – „se [[PERSON_33]]“, „o Soně Mikulkové“")
    ,@(27, "[[PERSON_34]] – „o [[PERSON_34]]“, „s [[PERSON_34]]“")
    ,@(30, "[[PERSON_35]] – „s [[PERSON_36]]“, „o [[PERSON_35]]“")
    ,@(31, "[[PERSON_37]] – „k [[PERSON_38]]“, „s [[PERSON_37]]“")
    ,@(32, "[[PERSON_39]] – „s [[PERSON_39]]“, „o [[PERSON_40]]“")
    ,@(33, "[[PERSON_41]] – „od [[PERSON_41]]“, „s [[PERSON_41]]“")
    ,@(34, "[[PERSON_42]] – „k [[PERSON_43]]“, „o [[PERSON_43]]“")
    ,@(35, "[[PERSON_44]] – „o [[PERSON_45]]“, „s [[PERSON_44]]“")
    ,@(36, "[[PERSON_46]] – „s [[PERSON_46]]“, „o [[PERSON_47]]“")
    ,@(37, "[[PERSON_48]] – „s [[PERSON_48]]“, „o [[PERSON_49]]“")
    ,@(38, "[[PERSON_50]] – „k [[PERSON_50]]“, „s [[PERSON_51]]“")
    ,@(39, "[[PERSON_52]] – „pro [[PERSON_53]]“, „o [[PERSON_54]]“")
    ,@(40, "[[PERSON_55]] – „k [[PERSON_55]]“, „o [[PERSON_55]]“")
    ,@(41, "[[PERSON_56]] – „o [[PERSON_57]]“, „s [[PERSON_56]]“")
    ,@(42, "[[PERSON_58]] – „s [[PERSON_59]]“, „o [[PERSON_60]]“")
    ,@(43, "[[PERSON_61]] – „s [[PERSON_61]]“, „o [[PERSON_61]]“")
    ,@(44, "[[PERSON_62]] – „u [[PERSON_62]]“, „o [[PERSON_63]]“")
    ,@(45, "[[PERSON_64]] – „se [[PERSON_65]]“, „o [[PERSON_64]]“")
    ,@(46, "[[PERSON_66]] – „o [[PERSON_67]]“, „s [[PERSON_68]]“")
    ,@(47, "[[PERSON_69]] – „k [[PERSON_70]]“, „o [[PERSON_70]]“")
    ,@(48, "[[PERSON_71]] – „o [[PERSON_72]]“, „s [[PERSON_71]]“")
    ,@(49, "[[PERSON_73]] – „s [[PERSON_73]]“, „o [[PERSON_73]]“")
    ,@(50, "[[PERSON_74]] – „s [[PERSON_74]]“, „o [[PERSON_75]]“")
    ,@(51, "[[PERSON_76]] – „o [[PERSON_77]]“, „s [[PERSON_78]]“")
    ,@(52, "[[PERSON_79]] – „s [[PERSON_79]]“, „o [[PERSON_80]]“")
    ,@(53, "[[PERSON_81]] – „o [[PERSON_81]]“, „s [[PERSON_81]]“")
    ,@(54, "[[PERSON_82]] – „s [[PERSON_82]]“, „o [[PERSON_83]]“")
)

foreach ($u in $updates) {
    $idx = $u[0]
    $newText = $u[1]
    $d.Paragraphs($idx).Range.Text = $newText
}

Write-Output "done"
